$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Add the new "step-2-see-url-title" table column after "step-2-see-url" ---
$newCol = $tbl.ListColumns.Add()
$newCol.Range.Cells.Item(1).Value = "step-2-see-url-title"

# --- Update the existing "step-2-see-url" column data: was numeric 4024, now the real URL text ---
$urlCol = $tbl.ListColumns.Item(19)
$urlCol.DataBodyRange.Value = "http://www.pbs.gov.au/info/industry/listing/elements/pbac-meetings/pbac-consumer-comments"

# --- Fill the new column's data body (after the URL column so shared-string order matches) ---
$newCol.DataBodyRange.Value = "PBAC Consumer Comment"

# --- Match formatting: header cell like the other header cells, data cells like other centered cells ---
$ws.Range("T1").Copy()
$newCol.Range.Cells.Item(1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("L3").Copy()
$newCol.DataBodyRange.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Copy formatting for the new column down into the blank rows below/above the table ---
$ws.Range("T11:T36").Copy()
$ws.Range("U11:U36").PasteSpecial(-4122)
$ws.Range("T44").Copy()
$ws.Range("U44").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row height tweaks (auto height grew once the new column/content was added) ---
$ws.Rows.Item(4).RowHeight = 105
$ws.Rows.Item(10).RowHeight = 105

# --- Selection matches where the edit was made ---
$ws.Range("U3:U10").Select()
